$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.579.58"
$ws.Range("E2").Value = "  +4.52%  "
$ws.Range("D3").Value = "3.393.36"
$ws.Range("E3").Value = "  +3.04%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'594.43"
$ws.Range("E5").Value = "  +7.59%  "
$ws.Range("D6").Value = "'186.63"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("D7").Value = "'0.601"
$ws.Range("E7").Value = "  +4.63%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "'0.184"
$ws.Range("E9").Value = "  +5.19%  "
$ws.Range("D10").Value = "'0.590"
$ws.Range("E10").Value = "  +2.66%  "
$ws.Range("D11").Value = "'47.60"
$ws.Range("E11").Value = "  +4.57%  "
$ws.Range("D12").Value = "'0.0000280"
$ws.Range("E12").Value = "  +7.49%  "
$ws.Range("D13").Value = "3.950.29"
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("D14").Value = "'638.73"
$ws.Range("E14").Value = "  +10.73%  "
$ws.Range("D15").Value = "'8.62"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "68.763.75"
$ws.Range("E16").Value = "  +4.73%  "
$ws.Range("D17").Value = "3.411.91"
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("D19").Value = "'18.08"
$ws.Range("E19").Value = "  +2.54%  "
$ws.Range("D20").Value = "'11.13"
$ws.Range("E20").Value = "  +3.18%  "
$ws.Range("D21").Value = "'0.913"
$ws.Range("E21").Value = "  +3.02%  "
$ws.Range("D22").Value = "'18.09"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").Value = "'5.10"
$ws.Range("E23").Value = "  +2.96%  "
$ws.Range("D24").Value = "'100.16"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("D25").Value = "'4.10"
$ws.Range("E25").Value = "  +4.86%  "
$ws.Range("D26").Value = "'2.86"
$ws.Range("E26").Value = "  +7.28%  "
$ws.Range("D27").Value = "'9.82"
$ws.Range("E27").Value = "  +5.42%  "
$ws.Range("D28").Value = "'32.85"
$ws.Range("E28").Value = "  +8.30%  "
$ws.Range("D29").Value = "'8.72"
$ws.Range("E29").Value = "  +4.66%  "
$ws.Range("D30").Value = "'6.87"
$ws.Range("E30").Value = "  +4.48%  "
$ws.Range("D31").Value = "'613.24"
$ws.Range("E31").Value = "  +8.44%  "
$ws.Range("D32").Value = "'3.81"
$ws.Range("E32").Value = "  +4.16%  "
$ws.Range("D33").Value = "4.033.67"
$ws.Range("E33").Value = "  +8.79%  "
$ws.Range("D34").Value = "'11.14"
$ws.Range("E34").Value = "  +3.19%  "
$ws.Range("D35").Value = "'0.106"
$ws.Range("E35").Value = "  +3.88%  "
$ws.Range("D36").Value = "'0.997"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "'56.79"
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("D38").Value = "'2.81"
$ws.Range("E38").Value = "  +9.17%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'3.32"
$ws.Range("E39").Value = "  +6.72%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.131"
$ws.Range("E40").Value = "  +4.87%  "
$ws.Range("D41").Value = "'33.75"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").Value = "0.0₃0709"
$ws.Range("E42").Value = "  +4.19%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "'3.43"
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.346"
$ws.Range("E44").Value = "  +4.38%  "
$ws.Range("D45").Value = "'0.0426"
$ws.Range("E45").Value = "  +5.07%  "
$ws.Range("D46").Value = "'0.130"
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("D47").Value = "'2.62"
$ws.Range("E47").Value = "  +5.05%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'1.39"
$ws.Range("E49").Value = "  +12.87%  "
$ws.Range("D50").Value = "'129.98"
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("D51").Value = "'7.82"
$ws.Range("E51").Value = "  +7.99%  "
